$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "71-18="
$t.Cell(1, 2).Range.Text = "3+83="
$t.Cell(1, 3).Range.Text = "71-44="
$t.Cell(1, 4).Range.Text = "28-17="
$t.Cell(1, 5).Range.Text = "95-10="
$t.Cell(2, 1).Range.Text = "7+76="
$t.Cell(2, 2).Range.Text = "99-47="
$t.Cell(2, 3).Range.Text = "91-2="
$t.Cell(2, 4).Range.Text = "1+19="
$t.Cell(2, 5).Range.Text = "14+16="
$t.Cell(3, 1).Range.Text = "97-11="
$t.Cell(3, 2).Range.Text = "95-23="
$t.Cell(3, 3).Range.Text = "91+0="
$t.Cell(3, 4).Range.Text = "17+77="
$t.Cell(3, 5).Range.Text = "21-9="
$t.Cell(4, 1).Range.Text = "51-34="
$t.Cell(4, 2).Range.Text = "96-77="
$t.Cell(4, 3).Range.Text = "79+11="
$t.Cell(4, 4).Range.Text = "94-59="
$t.Cell(4, 5).Range.Text = "94-4="
$t.Cell(5, 1).Range.Text = "62-43="
$t.Cell(5, 2).Range.Text = "65-64="
$t.Cell(5, 3).Range.Text = "3+31="
$t.Cell(5, 4).Range.Text = "87-10="
$t.Cell(5, 5).Range.Text = "75-31="
$t.Cell(6, 1).Range.Text = "38+49="
$t.Cell(6, 2).Range.Text = "94-8="
$t.Cell(6, 3).Range.Text = "85-24="
$t.Cell(6, 4).Range.Text = "94-89="
$t.Cell(6, 5).Range.Text = "22+56="
$t.Cell(7, 1).Range.Text = "68-50="
$t.Cell(7, 2).Range.Text = "74-71="
$t.Cell(7, 3).Range.Text = "55-10="
$t.Cell(7, 4).Range.Text = "18+10="
$t.Cell(7, 5).Range.Text = "69-23="
$t.Cell(8, 1).Range.Text = "75-27="
$t.Cell(8, 2).Range.Text = "26+19="
$t.Cell(8, 3).Range.Text = "31-8="
$t.Cell(8, 4).Range.Text = "73-68="
$t.Cell(8, 5).Range.Text = "83-23="
$t.Cell(9, 1).Range.Text = "34+48="
$t.Cell(9, 2).Range.Text = "92-69="
$t.Cell(9, 3).Range.Text = "66+9="
$t.Cell(9, 4).Range.Text = "82-78="
$t.Cell(9, 5).Range.Text = "66-20="
$t.Cell(10, 1).Range.Text = "21+54="
$t.Cell(10, 2).Range.Text = "73-67="
$t.Cell(10, 3).Range.Text = "69-68="
$t.Cell(10, 4).Range.Text = "75-47="
$t.Cell(10, 5).Range.Text = "98-57="
$t.Cell(11, 1).Range.Text = "83-65="
$t.Cell(11, 2).Range.Text = "1+81="
$t.Cell(11, 3).Range.Text = "71-47="
$t.Cell(11, 4).Range.Text = "85-84="
$t.Cell(11, 5).Range.Text = "27-26="
$t.Cell(12, 1).Range.Text = "52-35="
$t.Cell(12, 2).Range.Text = "75-53="
$t.Cell(12, 3).Range.Text = "91-50="
$t.Cell(12, 4).Range.Text = "0+79="
$t.Cell(12, 5).Range.Text = "50+2="
$t.Cell(13, 1).Range.Text = "33+60="
$t.Cell(13, 2).Range.Text = "93-36="
$t.Cell(13, 3).Range.Text = "42+39="
$t.Cell(13, 4).Range.Text = "88+3="
$t.Cell(13, 5).Range.Text = "11+2="
$t.Cell(14, 1).Range.Text = "34-29="
$t.Cell(14, 2).Range.Text = "33+6="
$t.Cell(14, 3).Range.Text = "44+42="
$t.Cell(14, 4).Range.Text = "82-32="
$t.Cell(14, 5).Range.Text = "9+31="
$t.Cell(15, 1).Range.Text = "3+56="
$t.Cell(15, 2).Range.Text = "62-48="
$t.Cell(15, 3).Range.Text = "72-37="
$t.Cell(15, 4).Range.Text = "28+39="
$t.Cell(15, 5).Range.Text = "80-17="
$t.Cell(16, 1).Range.Text = "17+5="
$t.Cell(16, 2).Range.Text = "21+74="
$t.Cell(16, 3).Range.Text = "97-76="
$t.Cell(16, 4).Range.Text = "4+31="
$t.Cell(16, 5).Range.Text = "14+5="
$t.Cell(17, 1).Range.Text = "84-28="
$t.Cell(17, 2).Range.Text = "45+18="
$t.Cell(17, 3).Range.Text = "24+29="
$t.Cell(17, 4).Range.Text = "64-40="
$t.Cell(17, 5).Range.Text = "52-25="
$t.Cell(18, 1).Range.Text = "41-40="
$t.Cell(18, 2).Range.Text = "0+81="
$t.Cell(18, 3).Range.Text = "85-70="
$t.Cell(18, 4).Range.Text = "56+14="
$t.Cell(18, 5).Range.Text = "42+11="
$t.Cell(19, 1).Range.Text = "14+46="
$t.Cell(19, 2).Range.Text = "19+46="
$t.Cell(19, 3).Range.Text = "45+14="
$t.Cell(19, 4).Range.Text = "89-1="
$t.Cell(19, 5).Range.Text = "33+51="
$t.Cell(20, 1).Range.Text = "1+17="
$t.Cell(20, 2).Range.Text = "43-9="
$t.Cell(20, 3).Range.Text = "29-24="
$t.Cell(20, 4).Range.Text = "13+32="
$t.Cell(20, 5).Range.Text = "23-0="
